$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C42").Font.Size = [double]8.5
Write-Host ("C42 font size=" + [string]$ws.Range("C42").Font.Size)
$ws.Range("E1").Font.Size = 10.5
Write-Host ("E1 font size=" + [string]$ws.Range("E1").Font.Size)
